# Update cryptocurrency price (D) and volume-change (E) columns
# with refreshed figures, matching a new data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are price strings formatted like "27.257.40" or "0.5213".
# Plain assignment lets Excel reinterpret single-dot strings as numbers
# (dropping significant trailing/leading zeros), so we force the cell to
# Text format before writing, then restore the default "Normal" style so
# no stray formatting is left behind.
$priceUpdates = @{
    'D2' = '27.257.40'
    'D3' = '1.898.61'
    'D4' = '1.003'
    'D5' = '307.94'
    'D7' = '0.5213'
    'D8' = '0.3774'
    'D9' = '0.07286'
    'D11' = '0.9019'
    'D12' = '0.08201'
    'D13' = '96.80'
    'D14' = '1.906.42'
    'D15' = '5.298'
    'D16' = '1.002'
    'D17' = '0.000008607'
    'D18' = '14.57'
    'D20' = '27.291.88'
    'D21' = '5.096'
    'D22' = '10.72'
    'D23' = '6.416'
    'D24' = '2.304'
    'D25' = '147.34'
    'D27' = '1.747'
    'D28' = '115.61'
    'D29' = '4.831'
    'D30' = '4.917'
    'D31' = '0.09245'
    'D32' = '0.05061'
    'D33' = '0.7972'
    'D34' = '1.233'
    'D35' = '3.436'
    'D36' = '2.967'
    'D37' = '2.593'
    'D38' = '0.5687'
    'D39' = '0.01998'
    'D41' = '8.981'
    'D42' = '6.572'
    'D43' = '115.37'
    'D44' = '0.1518'
    'D45' = '0.4887'
    'D47' = '10.06'
    'D48' = '1.624'
    'D49' = '38.16'
    'D50' = '63.68'
    'D51' = '0.05947'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Column E values are already non-numeric text (leading/trailing spaces,
# "+"/"-" sign, "%" suffix) so a direct assignment keeps them as strings.
$volumeUpdates = @{
    'E2' = '  +0.63%  '
    'E3' = '  +0.36%  '
    'E4' = '  +0.26%  '
    'E5' = '  +0.42%  '
    'E6' = '  +0.16%  '
    'E7' = '  +0.74%  '
    'E8' = '  +0.51%  '
    'E9' = '  +1.03%  '
    'E10' = '  +0.52%  '
    'E11' = '  +0.47%  '
    'E12' = '  +7.20%  '
    'E13' = '  +2.62%  '
    'E14' = '  +0.85%  '
    'E15' = '  +1.31%  '
    'E16' = '  +0.26%  '
    'E17' = '  +1.12%  '
    'E18' = '  +1.24%  '
    'E19' = '  +0.23%  '
    'E21' = '  +0.85%  '
    'E22' = '  +1.20%  '
    'E23' = '  +0.13%  '
    'E24' = '  +0.90%  '
    'E25' = '  +0.92%  '
    'E27' = '  +0.86%  '
    'E28' = '  +1.11%  '
    'E30' = '  -0.94%  '
    'E31' = '  +0.55%  '
    'E32' = '  +0.34%  '
    'E33' = '  +2.96%  '
    'E34' = '  -0.38%  '
    'E35' = '  +4.86%  '
    'E36' = '  -0.36%  '
    'E37' = '  +0.12%  '
    'E38' = '  +1.52%  '
    'E39' = '  +0.59%  '
    'E40' = '  +0.30%  '
    'E41' = '  -0.19%  '
    'E42' = '  -1.01%  '
    'E43' = '  -3.34%  '
    'E45' = '  +1.31%  '
    'E46' = '  +0.20%  '
    'E47' = '  -1.03%  '
    'E48' = '  +1.93%  '
    'E49' = '  +2.04%  '
    'E50' = '  -0.16%  '
    'E51' = '  +0.41%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

